$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header additions: P1=14, Q1=15 with style matching existing header (s=1, bold)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# New columns P and Q for rows 2-25, value 2 in every row
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 2   # column P
    $ws.Cells.Item($r, 17).Value = 2   # column Q
}

# Swap values in columns I, K, M, O for rows 2-25
# I: 1 -> 2 ; K: 2 -> 1 ; M: 1 -> 2 ; O: 2 -> 1
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 15).Value = 1   # O
}
